# Remove the 'PDB molecule' column from the metadata input file.
#
# The source sheet has columns:
#   A: Narrative ID
#   B: Object name (Genome AMA feature set)
#   C: Feature ID
#   D: PDB molecule        <-- column being removed
#   E: PDB filename
#   F: Is model
#
# Deleting column D shifts "PDB filename" and "Is model" left to D and E,
# matching the target layout (A:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "PDB molecule" column; cells to the right shift left.
$ws.Range("D1").EntireColumn.Delete()

# Widen the now-shifted "Object name" (B) and "Feature ID" (C) columns to fit
# their contents, as happened when the column was removed in Excel.
$ws.Columns.Item(2).ColumnWidth = 22.1640625
$ws.Columns.Item(3).ColumnWidth = 17.33203125

# Restore the cell selection left by the editor after the change.
$ws.Range("D10").Select()
